# menambahkan nama foto di template excel
# Insert a new column "Nama Photo" before the existing "KODE TIANG" column (G),
# shifting KODE TIANG/STATUS/Tarif one column to the right (G->H, H->I, I->J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at G; existing G,H,I (KODE TIANG, STATUS, Tarif) shift to H,I,J.
$ws.Columns("G").Insert()

# Header for the newly inserted column.
$ws.Range("G1").Value = "Nama Photo"

# Match the author's final selection in the saved workbook.
[void]$ws.Range("J10").Select()
